# Add "2022-Q3" quarter data:
#  - insert a new "2022-Q3" worksheet (fund-level holdings) right after the
#    "总计" summary sheet, pushing every other quarter sheet one slot later
#  - insert a matching summary row at the top of the "总计" sheet's data table

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert the new quarter as the first data row and
#    renumber the existing rows' index column (they keep their data, just
#    shift down by one row).
# ---------------------------------------------------------------------------

$summary = $wb.Worksheets.Item(1)

# Make room for a 9th data row by duplicating the formatting of the last
# existing row (this keeps column A's style, which carries the border/bold
# used for the index column, without disturbing anything else).
$summary.Range("A8:D8").Copy($summary.Range("A9:D9"))

# Now (re)write every data row top to bottom with the final values. Rows
# 2-8 already exist (only their values change); row 9 was just created
# above.
$summaryRows = @(
    @(0, "2022-Q3", 9,  5.61),
    @(1, "2022-Q2", 10, 5.44),
    @(2, "2022-Q1", 8,  5.21),
    @(3, "2021-Q4", 6,  5.34),
    @(4, "2021-Q3", 15, 6.21),
    @(5, "2021-Q2", 1,  5.35),
    @(6, "2021-Q1", 4,  5.5),
    @(7, "2020-Q4", 1,  7.12)
)

$r = 2
foreach ($row in $summaryRows) {
    $summary.Range("A$r").Value = $row[0]
    $summary.Range("B$r").Value = $row[1]
    $summary.Range("C$r").Value = $row[2]
    $summary.Range("D$r").Value = $row[3]
    $r++
}

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet with the fund-level holdings detail,
#    positioned right after "总计" (i.e. before the sheet that is currently
#    the "2022-Q2" tab).
# ---------------------------------------------------------------------------

$beforeSheet = $wb.Worksheets.Item(2)

$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Match the page-margin convention used by every other sheet in this
# workbook (0.75in/1in/0.5in) instead of the engine's brand-new-sheet
# defaults (0.7in/0.75in/0.3in).
$ps = $q3.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# Fetch the template (existing fund-detail sheet to copy formats from) by
# name *after* inserting the new tab, since worksheet lookup by numeric
# Item() index is position-based and would otherwise resolve to the
# sheet we just inserted.
$template = $wb.Worksheets.Item("2022-Q2")

# Copy the header row (labels + bold/border style) from the template sheet.
$template.Range("B1:H1").Copy($q3.Range("B1:H1"))

# Copy the index-column style (A2, bold border style "A" column uses) down
# for every data row we are about to fill in.
for ($i = 2; $i -le 10; $i++) {
    $template.Range("A2").Copy($q3.Range("A$i"))
}

# Columns B-G hold text values in this workbook (fund code / name / the
# numeric-looking figures are all stored as text), while column H (rank) is
# a real number. Force B:G to Text format *before* writing so values like
# "63.53" are not reinterpreted as numbers.
$q3.Range("B2:G10").NumberFormat = "@"

$q3Rows = @(
    @(0, "510810", "汇添富中证上海国企ETF",       "63.53", "98.46", "8.56", "5.4382", 2),
    @(1, "009073", "德邦惠利混合A",               "1.25",  "50.27", "2.91", "0.0364", 6),
    @(2, "001413", "中融鑫起点灵活配置混合A",      "0.59",  "77.02", "5.57", "0.0329", 4),
    @(3, "001739", "中融融安二号灵活配置混合",     "0.77",  "76.65", "3.96", "0.0305", 9),
    @(4, "015061", "中信建投沪深300指数增强A",     "1.47",  "91.17", "1.67", "0.0245", 5),
    @(5, "015062", "中信建投沪深300指数增强C",     "1.14",  "91.17", "1.67", "0.0190", 5),
    @(6, "001414", "中融鑫起点灵活配置混合C",      "0.21",  "77.02", "5.57", "0.0117", 4),
    @(7, "009074", "德邦惠利混合C",               "0.33",  "50.27", "2.91", "0.0096", 6),
    @(8, "010404", "博道盛利6个月持有期混合",      "1.10",  "41.15", "0.38", "0.0042", 9)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Range("A$r").Value = $row[0]
    $q3.Range("B$r").Value = $row[1]
    $q3.Range("C$r").Value = $row[2]
    $q3.Range("D$r").Value = $row[3]
    $q3.Range("E$r").Value = $row[4]
    $q3.Range("F$r").Value = $row[5]
    $q3.Range("G$r").Value = $row[6]
    $q3.Range("H$r").Value = $row[7]
    $r++
}
